# add server based classifier
#
# 1) Footer/date placeholder on the slide master and every slide layout:
#    "2020/7/20" -> "2020/8/9"
# 2) Slide 2 text box: "BessController" -> "SFFController" (and the
#    auto-fit textbox shrinks to match the new, shorter caption)
# 3) Slide 3 text box: "DockerController" -> "vnfController" (same
#    auto-fit shrink behaviour)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text everywhere it
#    appears: the slide master itself plus every one of its custom
#    layouts.
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2020/7/20") {
                $len = $tr.Text.Length
                $chars = $tr.Characters(1, $len)
                $chars.Text = "2020/8/9"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) & 3) Rename the two controller text boxes. Both shapes use
#    wrap="none" + spAutoFit, so PowerPoint itself would shrink the
#    box width to fit the new (shorter) caption; we reproduce the
#    measured target widths explicitly since this host does not run
#    real text layout. Left/Top/Height are untouched.
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $sh = $slide.Shapes.Item($j)
        if (-not $sh.HasTextFrame) {
            continue
        }
        $txt = $sh.TextFrame.TextRange.Text
        if ($txt -eq "BessController") {
            $sh.TextFrame.TextRange.Text = "SFFController"
            $sh.Width = 125.3624409448819   # 1592103 EMU
        } elseif ($txt -eq "DockerController") {
            $sh.TextFrame.TextRange.Text = "vnfController"
            $sh.Width = 123.4691353582677   # 1568058 EMU
        }
    }
}
